# "added leet code two sums problem"
# The sheet held a long list of product URLs (rows 1-29). The edit trims it
# down to a short "Two Sum"-style list of 9 entries: most rows are removed,
# a handful of the original rows are kept (their A-index renumbered to stay
# sequential), and a brand-new URL is appended as the final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove every row except the header (row 1) and the rows whose B-value
# survives into the final sheet (old rows 2,3,4,9,10,13,15,19). Delete from
# the bottom up so row numbers above the delete point stay stable.
$rowsToDelete = @(29,28,27,26,25,24,23,22,21,20,18,17,16,14,12,11,8,7,6,5)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Renumber column A (rows 2-9) back to a contiguous 0-based sequence.
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7

# Append the new row (index 8) with the new product URL.
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "https://www.uyyaala.com/products/enfamil-neuro-pro-gentlease-infant-formula-tin-pack-0-12m-777g"

# Match the formatting of the other index cells (bold, bordered, centered)
# by copying the format from the row above instead of rebuilding it
# property-by-property (which would create a new, near-duplicate style).
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122) | Out-Null

# Match the new selection recorded in the saved workbook.
$ws.Range("C1").Select() | Out-Null
